# #5: cash & deposit done
# Rework the "存款" (deposit) sheet: fix up the header row and append the
# same property_category/category/date/legislator_*/source_file/index
# metadata columns (G:M) that the other sheets (土地/建物/股票/事業投資)
# already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ---------------------------------------------------------------------
# Row 1: this row used to be a stray duplicate of row 2's literal values
# (B1="臺灣銀行東桃園分行", C1="活期儲蓄存款", D1="新臺幣", E1="廖正井",
# F1=2897128). Turn it into a proper header row, matching the other
# sheets' column-name convention, and extend it through column M.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 2).Value = "bank"
$ws.Cells.Item(1, 3).Value = "deposit_type"
$ws.Cells.Item(1, 4).Value = "currency"
$ws.Cells.Item(1, 5).Value = "owner"
$ws.Cells.Item(1, 6).Value = "total"
$ws.Cells.Item(1, 7).Value = "property_category"
$ws.Cells.Item(1, 8).Value = "category"
$ws.Cells.Item(1, 9).Value = "date"
$ws.Cells.Item(1, 10).Value = "legislator_name"
$ws.Cells.Item(1, 11).Value = "legislator_id"
$ws.Cells.Item(1, 12).Value = "source_file"
$ws.Cells.Item(1, 13).Value = "index"

# Copy the bold/centered header formatting (already present on B1:F1) onto
# the newly added G1:M1 header cells.
$ws.Cells.Item(1, 2).Copy() | Out-Null
$ws.Range("G1:M1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Data rows 2:25 -- bank / deposit_type / currency / owner / total were
# already correct; just make sure F11 is a real number (it used to be
# the text string "1041062") and then fill in the new metadata columns
# G:M for every row.
# ---------------------------------------------------------------------
$ws.Cells.Item(11, 6).Value = 1041062

$rows = @(
    @(2, 49, "臺灣銀行東桃園分行", "活期儲蓄存款", "廖正井", 2897128),
    @(3, 50, "臺灣土地銀行大安分行", "活期儲蓄存款", "廖正井", 1465357),
    @(4, 51, "第一商業銀行民生分行", "活期儲蓄存款", "廖正井", 1060910),
    @(5, 52, "台北富邦商業銀行市府分行", "活期儲蓄存款", "廖正井", 608422),
    @(6, 53, "台北富邦商業銀行城中分行", "活期儲蓄存款", "廖正井", 1366129),
    @(7, 54, "中華郵政股份有限公司桃圜府前郵局", "活期存款", "廖正井", 1466881),
    @(8, 55, "中華郵政股份有限公司台北信維郵局", "活期存款", "楊啟津", 269111),
    @(9, 56, "中華郵政股份有限公司台北信維郵局", "定期存款", "楊啟津", 1040087),
    @(10, 57, "國泰世華商業銀行信義分行", "定期存款", "楊啟津", 1136533),
    @(11, 58, "國泰世華商業銀行信義分行", "定期存款", "楊啟津", 1041062),
    @(12, 59, "國泰世華商業銀行信義分行", "活期儲蓄存款", "楊啟津", 1008078),
    @(13, 60, "臺灣土地銀行大安分行", "活期儲蓄存款", "楊啟津", 231391),
    @(14, 61, "臺灣土地銀行古亭分行", "活期儲蓄存款", "楊啟津", 227330),
    @(15, 62, "臺灣土地銀行大安分行", "定期存款", "楊啟津", 2822314),
    @(16, 64, "中國信託商業銀行敦南分行", "定期存款", "楊啟津", 3393392),
    @(17, 65, "台北富邦商業銀行敦和分行", "定期存款", "楊啟津", 1112015),
    @(18, 66, "兆豐國際商業銀行國外部", "活期儲蓄存款", "楊啟津", 1288752),
    @(19, 67, "華南商業銀行信維分行", "活期存款", "楊啟津", 7035),
    @(20, 68, "臺灣銀行大安分行", "活期存款", "楊啟津", 229407),
    @(21, 69, "遠東國際商業銀行營業部", "活期存款", "楊啟津", 717583),
    @(22, 70, "元大商業銀行南東分行", "活期存款", "楊啟津", 916972),
    @(23, 71, "臺灣銀行群賢分行", "活期存款", "廖正井", 138941),
    @(24, 72, "華南商業銀行信維分行", "定期存款", "楊啟津", 1094552),
    @(25, 73, "臺灣銀行大安分行", "定期存款", "楊啟津", 4000000)
)

foreach ($row in $rows) {
    $r = $row[0]

    # Columns B:F are unchanged from before -- re-assert them so the sheet
    # is self-consistent even if something upstream of this script nudged
    # them.
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = "新臺幣"
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]

    # New metadata columns, mirroring 土地/建物/股票/事業投資.
    $ws.Cells.Item($r, 7).Value = "deposit"
    $ws.Cells.Item($r, 8).Value = "normal"
    $ws.Cells.Item($r, 9).Value = "2013-12-20"
    $ws.Cells.Item($r, 10).Value = "廖正井"
    $ws.Cells.Item($r, 11).Value = 1711
    $ws.Cells.Item($r, 12).Value = "tmp393a1"
    $ws.Cells.Item($r, 13).Value = $row[1]
}

# Copy the plain data-row formatting (already present on B2:F25) onto the
# new G2:M25 block.
$ws.Cells.Item(2, 2).Copy() | Out-Null
$ws.Range("G2:M25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A1").Select() | Out-Null
